$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each string is one data row (columns A-T), joined by "|" to avoid comma/locale issues
$rowTexts = @(
  "FAPs|Agt|Mas1|ECs|3|1|5.788972999999999|17.366919|0.9474260381515466|0.9474260381515467|1|0.3333333333333333|0.1444666666666667|0.4334|0.06936556436665031|0.06936556436665031|0.8363136327333333|7.5268226946|0.0657187418320416|0.0657187418320416",
  "FAPs|Agt|Mas1|FAPs|3|1|5.788972999999999|17.366919|0.9474260381515466|0.9474260381515467|3|1|0.9857943333333333|2.957383|0.4733284283418029|0.4733284283418029|5.706736779219666|51.360631012977|0.4484436776083726|0.4484436776083726",
  "FAPs|Agt|Mas1|sCs|3|1|5.788972999999999|17.366919|0.9474260381515466|0.9474260381515467|3|1|0.9524246666666666|2.857274|0.4573060072915468|0.4573060072915468|5.513560679867332|49.622046118806|0.4332636187111324|0.4332636187111325",
  "sCs|Agt|Mas1|ECs|3|1|0.321238|0.963714|0.05257396184845335|0.05257396184845335|1|0.3333333333333333|0.1444666666666667|0.4334|0.06936556436665031|0.06936556436665031|0.04640818306666666|0.4176736476|0.003646822534608709|0.003646822534608709",
  "sCs|Agt|Mas1|FAPs|3|1|0.321238|0.963714|0.05257396184845335|0.05257396184845335|3|1|0.9857943333333333|2.957383|0.4733284283418029|0.4733284283418029|0.3166746000513333|2.850071400462|0.02488475073343033|0.02488475073343033",
  "sCs|Agt|Mas1|sCs|3|1|0.321238|0.963714|0.05257396184845335|0.05257396184845335|3|1|0.9524246666666666|2.857274|0.4573060072915468|0.4573060072915468|0.3059549950706666|2.753594955636|0.02404238858041431|0.02404238858041431"
)

$rowIdx = 2
foreach ($rowText in $rowTexts) {
    $parts = $rowText.Split("|")
    $colIdx = 1
    foreach ($p in $parts) {
        if ($colIdx -le 4) {
            $ws.Cells.Item($rowIdx, $colIdx).Value = $p
        } else {
            $ws.Cells.Item($rowIdx, $colIdx).Value = [double]$p
        }
        $colIdx++
    }
    $rowIdx++
}
